$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(987.7962036132812, 0.741, 0.7583000063896179, 1, 0.3404000103473663)
    3 = @(1037.638916015625, 0.8300999999999999, 0.8327, 1, 0.6093000173568726)
    4 = @(659.8099975585938, 0.7708, 0.77, 0.9927999973297119, 0.5852000117301941)
    5 = @(609.0084838867188, 0.6424, 0.6657, 0.9564999938011169, 0.2389000058174133)
    6 = @(787.65380859375, 0.6231, 0.6476, 0.9884999990463257, 0.1343999952077866)
    7 = @(642.2716064453125, 0.6468, 0.6626999974250793, 0.8011999726295471, 0.4117999970912933)
    8 = @(668.5167846679688, 0.599, 0.5961, 0.7271999716758728, 0.4117999970912933)
    9 = @(5392.69580078125, 0.6951000000000001, 0.6861, 1, 0.1343999952077866)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
}
